# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" (rId1), "演出" (rId2) and "全部类型" (rId4) sheets to match the
# latest scrape output.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12781
$ws1.Range("F3").Value  = 620
$ws1.Range("F5").Value  = 29
$ws1.Range("F6").Value  = 314
$ws1.Range("F7").Value  = 401
$ws1.Range("F9").Value  = 12783
$ws1.Range("F11").Value = 16
$ws1.Range("F12").Value = 5206
$ws1.Range("F14").Value = 17
$ws1.Range("F15").Value = 11
$ws1.Range("F16").Value = 26
$ws1.Range("F17").Value = 1199
$ws1.Range("F18").Value = 30
$ws1.Range("F19").Value = 132
$ws1.Range("F20").Value = 669
$ws1.Range("F22").Value = 6156
$ws1.Range("F24").Value = 3620
$ws1.Range("F25").Value = 220
$ws1.Range("F26").Value = 40

# ---- Sheet: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 19

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 12781
$ws4.Range("F3").Value  = 620
$ws4.Range("F5").Value  = 29
$ws4.Range("F6").Value  = 314
$ws4.Range("F7").Value  = 19
$ws4.Range("F8").Value  = 401
$ws4.Range("F10").Value = 12783
$ws4.Range("F12").Value = 16
$ws4.Range("F13").Value = 5206
$ws4.Range("F15").Value = 17
$ws4.Range("F16").Value = 11
$ws4.Range("F17").Value = 26
$ws4.Range("F18").Value = 1199
$ws4.Range("F19").Value = 30
$ws4.Range("F20").Value = 132
$ws4.Range("F21").Value = 669
$ws4.Range("F24").Value = 6156
$ws4.Range("F26").Value = 3620
$ws4.Range("F27").Value = 220
$ws4.Range("F28").Value = 40
